$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 251, shifting existing rows 251-270 down to 252-271
$ws.Rows("251:251").Insert()

# Populate the new row 251 with the new weekly data entry
$ws.Cells.Item(251, 1).Value = 3
$ws.Cells.Item(251, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(251, 3).Value = "Coquimbo"
$ws.Cells.Item(251, 4).Value = 44461
$ws.Cells.Item(251, 5).Value = 5
$ws.Cells.Item(251, 6).Value = 100112021
$ws.Cells.Item(251, 7).Value = "Ají"
$ws.Cells.Item(251, 8).Value = "Americana (o)"
$ws.Cells.Item(251, 9).Value = "Primera"
$ws.Cells.Item(251, 10).Value = 35
$ws.Cells.Item(251, 11).Value = 41000
$ws.Cells.Item(251, 12).Value = 41000
$ws.Cells.Item(251, 13).Value = 41000
$ws.Cells.Item(251, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(251, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(251, 16).Value = 2733
$ws.Cells.Item(251, 17).Value = 15
$ws.Cells.Item(251, 18).Value = "Hortaliza"
